$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 343.17142
$ws.Range("J17").Value = 285.66666
$ws.Range("L17").Value = 856.9999799999999
$ws.Range("N17").Value = -1192.99998
$ws.Range("H96").Value = 1292.8667
$ws.Range("I96").Value = 854
$ws.Range("K96").Value = 2562
$ws.Range("M96").Value = -1189
$ws.Range("H97").Value = 898.7143
$ws.Range("J97").Value = 898.7143
$ws.Range("L97").Value = 2696.1429
$ws.Range("N97").Value = -3688.1429
$ws.Range("H98").Value = 721.5
$ws.Range("J98").Value = 994.3333
$ws.Range("L98").Value = 994.3333
$ws.Range("N98").Value = -3990.3333
$ws.Range("H106").Value = 2931.0476
$ws.Range("I106").Value = 3218
$ws.Range("K106").Value = 3218
$ws.Range("M106").Value = -2587
$ws.Range("H116").Value = 7518.85
$ws.Range("I116").Value = 5366.6665
$ws.Range("J116").Value = 8441.214
$ws.Range("K116").Value = 5366.6665
$ws.Range("L116").Value = 8441.214
$ws.Range("M116").Value = -1924.6665
$ws.Range("N116").Value = -15325.214
$ws.Range("H122").Value = 721.5
$ws.Range("J122").Value = 994.3333
$ws.Range("L122").Value = 2982.9999
$ws.Range("N122").Value = -7882.9999
$ws.Range("H131").Value = 10762
$ws.Range("I131").Value = 8764.615
$ws.Range("K131").Value = 26293.845
$ws.Range("M131").Value = -21253.845
$ws.Range("H132").Value = 11040.604
$ws.Range("I132").Value = 1284.7347
$ws.Range("K132").Value = 3854.2041
$ws.Range("M132").Value = -1324.2041
$ws.Range("H138").Value = 2527.76
$ws.Range("I138").Value = 1655.875
$ws.Range("J138").Value = 4077.7778
$ws.Range("K138").Value = 4967.625
$ws.Range("L138").Value = 12233.3334
$ws.Range("M138").Value = 172.375
$ws.Range("N138").Value = -22513.3334
$ws.Range("H139").Value = 44197.777
$ws.Range("J139").Value = 44197.777
$ws.Range("L139").Value = 44197.777
$ws.Range("N139").Value = -54477.777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 34682
$ws.Range("I43").Value = 19342
$ws.Range("J43").Value = 39795.332
$ws.Range("K43").Value = 19342
$ws.Range("L43").Value = 39795.332
$ws.Range("M43").Value = -19029
$ws.Range("N43").Value = -40421.332
$ws.Range("H45").Value = 5649.8184
$ws.Range("J45").Value = 10666.333
$ws.Range("L45").Value = 10666.333
$ws.Range("N45").Value = -11420.333
$ws.Range("H61").Value = 8059.579
$ws.Range("I61").Value = 1519
$ws.Range("K61").Value = 1519
$ws.Range("M61").Value = -1307
$ws.Range("H76").Value = 49969
$ws.Range("J76").Value = 49969
$ws.Range("L76").Value = 49969
$ws.Range("N76").Value = -50645
$ws.Range("H79").Value = 49969
$ws.Range("J79").Value = 49969
$ws.Range("L79").Value = 49969
$ws.Range("N79").Value = -52309
$ws.Range("H132").Value = 1241.2307
$ws.Range("I132").Value = 883.5862
$ws.Range("J132").Value = 2278.4
$ws.Range("K132").Value = 2650.7586
$ws.Range("L132").Value = 6835.200000000001
$ws.Range("M132").Value = -120.7586000000001
$ws.Range("N132").Value = -11895.2
$ws.Range("H136").Value = 8059.579
$ws.Range("I136").Value = 1519
$ws.Range("K136").Value = 4557
$ws.Range("M136").Value = -2007

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 18557224
$ws.Range("I86").Value = 35787184
$ws.Range("K86").Value = 35787184
$ws.Range("M86").Value = -35786061
$ws.Range("H88").Value = 37012.8
$ws.Range("J88").Value = 40766
$ws.Range("L88").Value = 40766
$ws.Range("N88").Value = -41578
$ws.Range("H89").Value = 18557224
$ws.Range("I89").Value = 35787184
$ws.Range("K89").Value = 178935920
$ws.Range("M89").Value = -178930304
$ws.Range("H91").Value = 37012.8
$ws.Range("J91").Value = 40766
$ws.Range("L91").Value = 40766
$ws.Range("N91").Value = -43574
$ws.Range("H105").Value = 2877.5715
$ws.Range("I105").Value = 2198
$ws.Range("J105").Value = 3387.25
$ws.Range("K105").Value = 2198
$ws.Range("L105").Value = 3387.25
$ws.Range("M105").Value = -451
$ws.Range("N105").Value = -6881.25
$ws.Range("H134").Value = 2195.8333
$ws.Range("I134").Value = 1795.1428
$ws.Range("K134").Value = 5385.428400000001
$ws.Range("M134").Value = -2850.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2883.7222
$ws.Range("I58").Value = 1910.8572
$ws.Range("K58").Value = 1910.8572
$ws.Range("M58").Value = -1707.8572
$ws.Range("H105").Value = 1200
$ws.Range("I105").Value = 1200
$ws.Range("J105").Value = 1200
$ws.Range("K105").Value = 1200
$ws.Range("L105").Value = 1200
$ws.Range("M105").Value = 547
$ws.Range("N105").Value = -4694
$ws.Range("H107").Value = 6890.0625
$ws.Range("I107").Value = 256.875
$ws.Range("J107").Value = 13523.25
$ws.Range("K107").Value = 256.875
$ws.Range("L107").Value = 13523.25
$ws.Range("M107").Value = 1663.125
$ws.Range("N107").Value = -17363.25
$ws.Range("H132").Value = 2803.2632
$ws.Range("I132").Value = 2648.1667
$ws.Range("J132").Value = 5595
$ws.Range("K132").Value = 7944.500100000001
$ws.Range("L132").Value = 16785
$ws.Range("M132").Value = -5414.500100000001
$ws.Range("N132").Value = -21845
$ws.Range("H134").Value = 6250.25
$ws.Range("I134").Value = 5562.85
$ws.Range("K134").Value = 16688.55
$ws.Range("M134").Value = -14153.55
$ws.Range("H136").Value = 2883.7222
$ws.Range("I136").Value = 1910.8572
$ws.Range("K136").Value = 5732.571599999999
$ws.Range("M136").Value = -3182.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 1000
$ws.Range("J52").Value = 1000
$ws.Range("L52").Value = 3000
$ws.Range("N52").Value = -3532

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 600
$ws.Range("J29").Value = 600
$ws.Range("L29").Value = 600
$ws.Range("N29").Value = -1180
$ws.Range("H41").Value = 10249.5
$ws.Range("I41").Value = 1500
$ws.Range("K41").Value = 1500
$ws.Range("M41").Value = -1145
$ws.Range("H70").Value = 192799.83
$ws.Range("I70").Value = 285975
$ws.Range("K70").Value = 285975
$ws.Range("M70").Value = -285705
$ws.Range("H73").Value = 192799.83
$ws.Range("I73").Value = 285975
$ws.Range("K73").Value = 285975
$ws.Range("M73").Value = -285039
$ws.Range("H132").Value = 5922.608
$ws.Range("I132").Value = 5376.8047
$ws.Range("J132").Value = 8160.4
$ws.Range("K132").Value = 16130.4141
$ws.Range("L132").Value = 24481.2
$ws.Range("M132").Value = -13600.4141
$ws.Range("N132").Value = -29541.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 10101829
$ws.Range("I22").Value = 12987638
$ws.Range("J22").Value = 1496
$ws.Range("K22").Value = 12987638
$ws.Range("L22").Value = 1496
$ws.Range("M22").Value = -12987343
$ws.Range("N22").Value = -2086
$ws.Range("H27").Value = 10101829
$ws.Range("I27").Value = 12987638
$ws.Range("J27").Value = 1496
$ws.Range("K27").Value = 12987638
$ws.Range("L27").Value = 1496
$ws.Range("M27").Value = -12987531
$ws.Range("N27").Value = -1710
$ws.Range("H46").Value = 2614.1428
$ws.Range("J46").Value = 3433.111
$ws.Range("L46").Value = 3433.111
$ws.Range("N46").Value = -3809.111
$ws.Range("H55").Value = 477.36365
$ws.Range("I55").Value = 433.42856
$ws.Range("J55").Value = 554.25
$ws.Range("K55").Value = 433.42856
$ws.Range("L55").Value = 554.25
$ws.Range("M55").Value = -260.42856
$ws.Range("N55").Value = -900.25
$ws.Range("H82").Value = 40001840
$ws.Range("I82").Value = 62501988
$ws.Range("K82").Value = 62501988
$ws.Range("M82").Value = -62501627
$ws.Range("H85").Value = 40001840
$ws.Range("I85").Value = 62501988
$ws.Range("K85").Value = 62501988
$ws.Range("M85").Value = -62500740
$ws.Range("H122").Value = 4118.0713
$ws.Range("I122").Value = 3263.75
$ws.Range("J122").Value = 5257.1665
$ws.Range("K122").Value = 9791.25
$ws.Range("L122").Value = 15771.4995
$ws.Range("M122").Value = -7341.25
$ws.Range("N122").Value = -20671.4995
$ws.Range("H136").Value = 2843.4614
$ws.Range("I136").Value = 2891
$ws.Range("K136").Value = 8673
$ws.Range("M136").Value = -6123

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 16000
$ws.Range("I21").Value = 16000
$ws.Range("K21").Value = 16000
$ws.Range("M21").Value = -15765
$ws.Range("H35").Value = 16000
$ws.Range("I35").Value = 16000
$ws.Range("K35").Value = 16000
$ws.Range("M35").Value = -15710
$ws.Range("H81").Value = 8549213
$ws.Range("I81").Value = 1684.3334
$ws.Range("J81").Value = 27781152
$ws.Range("K81").Value = 3368.6668
$ws.Range("L81").Value = 55562304
$ws.Range("M81").Value = -2307.6668
$ws.Range("N81").Value = -55564426
$ws.Range("H84").Value = 8549213
$ws.Range("I84").Value = 1684.3334
$ws.Range("J84").Value = 27781152
$ws.Range("K84").Value = 16843.334
$ws.Range("L84").Value = 277811520
$ws.Range("M84").Value = -11539.334
$ws.Range("N84").Value = -277822128
$ws.Range("H122").Value = 2586.25
$ws.Range("I122").Value = 1972.5
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 5917.5
$ws.Range("L122").Value = 9600
$ws.Range("M122").Value = -3467.5
$ws.Range("N122").Value = -14500
$ws.Range("H136").Value = 49734.7
$ws.Range("I136").Value = 65036.332
$ws.Range("J136").Value = 3829.8
$ws.Range("K136").Value = 195108.996
$ws.Range("L136").Value = 11489.4
$ws.Range("M136").Value = -192558.996
$ws.Range("N136").Value = -16589.4
